$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 30,8
$arr[0,0] = 0
$arr[0,1] = "walkingToRunning"
$arr[0,2] = 7.460176270583601
$arr[0,3] = 1.940670852003414
$arr[0,4] = -2.248238760849495
$arr[0,5] = -1.420729875564575
$arr[0,6] = 0.4241084158420563
$arr[0,7] = 0.8468852639198303

$arr[1,0] = 100
$arr[1,1] = "walkingToRunning"
$arr[1,2] = 7.05335301366344
$arr[1,3] = 1.372367842444078
$arr[1,4] = -2.086995552326103
$arr[1,5] = -0.8428239226341248
$arr[1,6] = 1.781787991523743
$arr[1,7] = 1.107874989509582

$arr[2,0] = 200
$arr[2,1] = "walkingToRunning"
$arr[2,2] = 10.34839492008603
$arr[2,3] = 0.2938980234079889
$arr[2,4] = 1.988064823479517
$arr[2,5] = 2.085587978363037
$arr[2,6] = -0.4573979377746582
$arr[2,7] = -1.845570683479309

$arr[3,0] = 300
$arr[3,1] = "walkingToRunning"
$arr[3,2] = 13.00025196733149
$arr[3,3] = 5.53139660276216
$arr[3,4] = -0.4643907300357014
$arr[3,5] = 1.892242550849914
$arr[3,6] = -2.681137561798096
$arr[3,7] = -5.550560474395752

$arr[4,0] = 400
$arr[4,1] = "walkingToRunning"
$arr[4,2] = 14.22385301261109
$arr[4,3] = 6.654247678559409
$arr[4,4] = -0.2390175523429294
$arr[4,5] = 2.531933784484864
$arr[4,6] = 1.24702525138855
$arr[4,7] = -2.802178144454956

$arr[5,0] = 500
$arr[5,1] = "walkingToRunning"
$arr[5,2] = 3.261832862064715
$arr[5,3] = 4.847148434868966
$arr[5,4] = -3.660793435984526
$arr[5,5] = 2.003562688827514
$arr[5,6] = 2.79219126701355
$arr[5,7] = -3.079679489135742

$arr[6,0] = 600
$arr[6,1] = "walkingToRunning"
$arr[6,2] = 2.012637302793288
$arr[6,3] = -7.111074842255716
$arr[6,4] = 3.683522890354033
$arr[6,5] = -6.388058185577393
$arr[6,6] = -2.879809379577637
$arr[6,7] = 3.848800420761109

$arr[7,0] = 700
$arr[7,1] = "walkingToRunning"
$arr[7,2] = 12.92323684692373
$arr[7,3] = -4.288529297400908
$arr[7,4] = -1.309652714893691
$arr[7,5] = -14.3317346572876
$arr[7,6] = -11.627760887146
$arr[7,7] = 0.5219795703887939

$arr[8,0] = 800
$arr[8,1] = "walkingToRunning"
$arr[8,2] = 33.17913206692398
$arr[8,3] = -49.48561506846823
$arr[8,4] = -8.165793155801595
$arr[8,5] = 6.343982696533203
$arr[8,6] = -4.255598545074463
$arr[8,7] = -2.953978300094604

$arr[9,0] = 900
$arr[9,1] = "walkingToRunning"
$arr[9,2] = 29.25954614836614
$arr[9,3] = -1.534980112108286
$arr[9,4] = 10.88253550693887
$arr[9,5] = 2.917559623718262
$arr[9,6] = -4.095808506011963
$arr[9,7] = -0.09321063756942739

$arr[10,0] = 1000
$arr[10,1] = "walkingToRunning"
$arr[10,2] = 8.380168257088386
$arr[10,3] = -41.39767669809267
$arr[10,4] = 28.80079953423849
$arr[10,5] = 7.193531036376953
$arr[10,6] = -2.798849105834961
$arr[10,7] = -2.648780107498169

$arr[11,0] = 1100
$arr[11,1] = "walkingToRunning"
$arr[11,2] = -35.82913970947266
$arr[11,3] = -7.401230812072754
$arr[11,4] = -22.76275253295898
$arr[11,5] = 0.068509817123413
$arr[11,6] = -2.387657165527344
$arr[11,7] = 1.848233819007873

$arr[12,0] = 1200
$arr[12,1] = "walkingToRunning"
$arr[12,2] = -0.386668369687925
$arr[12,3] = 6.335032199991126
$arr[12,4] = -3.633465668250821
$arr[12,5] = -5.729191780090332
$arr[12,6] = -0.3087935447692871
$arr[12,7] = 0.4756405651569366

$arr[13,0] = 1300
$arr[13,1] = "walkingToRunning"
$arr[13,2] = -0.332839686295094
$arr[13,3] = -6.914240393145466
$arr[13,4] = -1.636375197048832
$arr[13,5] = -3.537942886352539
$arr[13,6] = -2.56502366065979
$arr[13,7] = -3.619768619537354

$arr[14,0] = 1400
$arr[14,1] = "walkingToRunning"
$arr[14,2] = -1.074484512723721
$arr[14,3] = -6.283810664867543
$arr[14,4] = 30.66327963204198
$arr[14,5] = 3.089067220687866
$arr[14,6] = 14.31196117401123
$arr[14,7] = -0.1933455020189285

$arr[15,0] = 1500
$arr[15,1] = "walkingToRunning"
$arr[15,2] = 4.744677149016312
$arr[15,3] = 18.54307062872537
$arr[15,4] = 23.11713764585292
$arr[15,5] = 3.182277917861938
$arr[15,6] = -2.187387466430664
$arr[15,7] = -2.911900281906128

$arr[16,0] = 1600
$arr[16,1] = "walkingToRunning"
$arr[16,2] = 40.11636086167968
$arr[16,3] = -43.34975341270727
$arr[16,4] = 39.71557577725119
$arr[16,5] = 3.544467687606812
$arr[16,6] = -1.478453874588013
$arr[16,7] = -2.32813549041748

$arr[17,0] = 1700
$arr[17,1] = "walkingToRunning"
$arr[17,2] = 8.228727192714175
$arr[17,3] = 2.979075727791802
$arr[17,4] = -10.42038246680967
$arr[17,5] = -2.930742263793945
$arr[17,6] = 2.399108648300171
$arr[17,7] = 2.079396247863769

$arr[18,0] = 1800
$arr[18,1] = "walkingToRunning"
$arr[18,2] = 8.725867090554029
$arr[18,3] = 5.052738584321171
$arr[18,4] = -5.431183058640126
$arr[18,5] = -15.79913711547852
$arr[18,6] = -0.4046673476696014
$arr[18,7] = 7.474960803985596

$arr[19,0] = 1900
$arr[19,1] = "walkingToRunning"
$arr[19,2] = 36.49437555773487
$arr[19,3] = -39.19836583630763
$arr[19,4] = 36.29759532007652
$arr[19,5] = 6.207096099853516
$arr[19,6] = -14.30157470703125
$arr[19,7] = 1.141430854797363

$arr[20,0] = 2000
$arr[20,1] = "walkingToRunning"
$arr[20,2] = 17.42426753866236
$arr[20,3] = -11.81626648738428
$arr[20,4] = 23.43557239400905
$arr[20,5] = 2.070674419403076
$arr[20,6] = -3.818307399749756
$arr[20,7] = 2.332929134368896

$arr[21,0] = 2100
$arr[21,1] = "walkingToRunning"
$arr[21,2] = 13.20562122608059
$arr[21,3] = -5.217379109612901
$arr[21,4] = 31.13059931787963
$arr[21,5] = 12.43622970581055
$arr[21,6] = 1.804691195487976
$arr[21,7] = -0.571514368057251

$arr[22,0] = 2200
$arr[22,1] = "walkingToRunning"
$arr[22,2] = -9.915741953356811
$arr[22,3] = 7.114933737393288
$arr[22,4] = 0.5023386725059567
$arr[22,5] = -4.415188312530518
$arr[22,6] = 2.631336450576782
$arr[22,7] = -0.9896306991577148

$arr[23,0] = 2300
$arr[23,1] = "walkingToRunning"
$arr[23,2] = -8.703923554256075
$arr[23,3] = 12.09582072290882
$arr[23,4] = -13.1394268562054
$arr[23,5] = -4.100935459136963
$arr[23,6] = 1.459545493125916
$arr[23,7] = 7.553257465362549

$arr[24,0] = 2400
$arr[24,1] = "walkingToRunning"
$arr[24,2] = -1.12605124506454
$arr[24,3] = -8.641085657580049
$arr[24,4] = -4.591965708239329
$arr[24,5] = -3.94806981086731
$arr[24,6] = 8.382166862487793
$arr[24,7] = -2.697782278060913

$arr[25,0] = 2500
$arr[25,1] = "walkingToRunning"
$arr[25,2] = -3.980940161080194
$arr[25,3] = -14.64417253691573
$arr[25,4] = 2.824098981660002
$arr[25,5] = 4.418517112731934
$arr[25,6] = 15.49866580963135
$arr[25,7] = -7.229417324066162

$arr[26,0] = 2600
$arr[26,1] = "walkingToRunning"
$arr[26,2] = 7.344494918297197
$arr[26,3] = -4.105810165405226
$arr[26,4] = 11.71896579347811
$arr[26,5] = 1.811815142631531
$arr[26,6] = -5.760817050933838
$arr[26,7] = -5.106345176696777

$arr[27,0] = 2700
$arr[27,1] = "walkingToRunning"
$arr[27,2] = 22.1531091229668
$arr[27,3] = -10.78129972260548
$arr[27,4] = 11.93481448601038
$arr[27,5] = 6.468618869781494
$arr[27,6] = 1.199088335037231
$arr[27,7] = -6.144978046417236

$arr[28,0] = 2800
$arr[28,1] = "walkingToRunning"
$arr[28,2] = 16.40996650169633
$arr[28,3] = -8.967308833681265
$arr[28,4] = 6.306509708536076
$arr[28,5] = -2.402371168136597
$arr[28,6] = 2.291517019271851
$arr[28,7] = 1.192030906677246

$arr[29,0] = 2900
$arr[29,1] = "walkingToRunning"
$arr[29,2] = -1.831812381744486
$arr[29,3] = 3.093760391761203
$arr[29,4] = 3.286617870988547
$arr[29,5] = -12.0989408493042
$arr[29,6] = -13.58198833465576
$arr[29,7] = -2.806971788406372

$ws.Range("A2:H31").Value = $arr

Write-Output "done"